# 2024_movies.xlsx — "Add files via upload" re-upload.
#
# The re-uploaded workbook re-ran the same Power Query / data-prep pipeline
# that produced the sheet, which emitted the feature columns (C:M) in a
# different order than before (the identifier/date/name columns A,B and the
# numeric-ROI/revenue columns N,O keep their positions). Column C, which
# used to be a plain 0/1 "is sequel" flag, is now the "Month Name" text
# column; the former "Month Name" data (and the rest of C:M) shifts
# accordingly. The O-column "Predicted Revenue" formulas are retyped to
# reference the Production Budget column at its new address (L instead of
# I), and column N's width shrinks back to its default-ish size.
#
# Net effect implemented below: rewrite the header row and all data rows
# for columns C:M with their new contents, re-point the O formulas at L,
# resize column N, and move the active selection to O13 (matching the
# saved sheetView state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header text for C1:M1 (A1/B1/N1/O1 are unchanged).
$headers = @("Month Name", "Franchise", "Running Time", "Running Time Group", `
    "Genre", "Creative Type", "Distributor", "Distributor Type", "Mpaa", `
    "Production Budget", "Budget Range")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 3 + $i).Value = $headers[$i]
}

# New C:M content per data row (row number => values for columns C..M).
$rows = @{
    2 = @("nov", 0, 123, "100-125", "Action", "Contemporary Fiction", "Amazon", "Medium", "PG-13", 200000000, "200-225M")
    3 = @("may", 1, 145, "125-150", "Action", "Science Fiction", "Disney", "Big", "PG-13", 160000000, "150-175M")
    4 = @("mar", 1, 115, "100-125", "Action", "Science Fiction", "Warner Bros.", "Big", "PG-13", 135000000, "125-150M")
    5 = @("nov", 1, 160, "150-175", "Musical", "Fantasy", "Universal", "Big", "PG", 145000000, "125-150M")
    6 = @("oct", 0, 104, "100-125", "Drama", "Historical Fiction", "Sony", "Big", "PG-13", 50000000, "50-75M")
    7 = @("oct", 0, 120, "100-125", "Drama", "Dramatization", "Briarcliff Entertainment", "Small", "R", 16000000, "15-25M")
    8 = @("oct", 1, 109, "100-125", "Action", "Super Hero", "Sony", "Big", "PG-13", 110000000, "100-125M")
    9 = @("oct", 1, 138, "125-150", "Thriller/Suspense", "Contemporary Fiction", "Warner Bros.", "Big", "R", 190000000, "175-200M")
}

foreach ($r in 2..9) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 3 + $i).Value = $vals[$i]
    }
}

# O-column "Predicted Revenue" formulas now reference the Production
# Budget column at its new position (L) instead of the old one (I).
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 15).Formula = "=N$r*L$r + L$r"
}

# Column N (index 14) width shrinks from the old custom width down to
# (approximately) the default column width.
$ws.Columns.Item(14).ColumnWidth = 11.6666667

# Active cell / selection moves to O13.
$ws.Range("O13").Select()
